# Generate Report for Handoff
# Replace the old GUID-named source file references with the new GUID,
# and bump the handoff/generate timestamps, across the Overview, zh-cn and
# de-de sheets (mirrors the localization-status.xlsx handback report).

$wb = $excel.ActiveWorkbook

$oldGuidFile   = "7dfa5e28-132e-416c-bcdc-effcd08ed1b2.md"
$newGuidFile   = "521dd8e8-7173-4a8f-a108-8f5b025fee41.md"
$oldGuidFileRel = "e2e\7dfa5e28-132e-416c-bcdc-effcd08ed1b2.md"
$newGuidFileRel = "e2e\521dd8e8-7173-4a8f-a108-8f5b025fee41.md"

$hyperlinkAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee937ad28106c2e8463230d161049af3037dc28e/e2e/7dfa5e28-132e-416c-bcdc-effcd08ed1b2.md"

$oldZhXlf = "7dfa5e28-132e-416c-bcdc-effcd08ed1b2.d9e0bfdcfb351fc317470feabfc42e43cdd006fb.zh-cn.xlf"
$newZhXlf = "521dd8e8-7173-4a8f-a108-8f5b025fee41.d3b116b757a0e2eed52ac25728cb7fe2c40de290.zh-cn.xlf"
$oldDeXlf = "7dfa5e28-132e-416c-bcdc-effcd08ed1b2.d9e0bfdcfb351fc317470feabfc42e43cdd006fb.de-de.xlf"
$newDeXlf = "521dd8e8-7173-4a8f-a108-8f5b025fee41.d3b116b757a0e2eed52ac25728cb7fe2c40de290.de-de.xlf"

$oldGenerateDate = "2016-08-22 01:07:04"
$newGenerateDate = "2016-08-22 01:07:36"

$oldZhHandoffDate = "2016-08-22 01:06:56"
$newZhHandoffDate = "2016-08-22 01:07:26"

# ---------------- Overview sheet ----------------
$wsOverview = $wb.Worksheets.Item("Overview")

# A2: plain file name
$wsOverview.Range("A2").Value = $newGuidFile

# B2: "e2e\<guid>.md" -- also carries the external hyperlink, so rebuild the
# hyperlink (delete + re-add) so both the cell text and the hyperlink's
# display text move to the new guid while keeping the same target address.
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddr, "", "", $newGuidFileRel)

# G2: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = $newGenerateDate

# ---------------- zh-cn sheet ----------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# A2: hyperlink + display is the plain file name (no e2e\ prefix here)
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hyperlinkAddr, "", "", $newGuidFile)

# G2: Latest Handoff File
$wsZh.Range("G2").Value = $newZhXlf

# H2: Latest Handoff Datetime
$wsZh.Range("H2").Value = $newZhHandoffDate

# ---------------- de-de sheet ----------------
$wsDe = $wb.Worksheets.Item("de-de")

# A2: hyperlink + display is the plain file name (no e2e\ prefix here)
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hyperlinkAddr, "", "", $newGuidFile)

# G2: Latest Handoff File (de-de xlf)
$wsDe.Range("G2").Value = $newDeXlf

# H2: Latest Handback DateTime (shared with Overview's G2 generate date)
$wsDe.Range("H2").Value = $newGenerateDate
